# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210
#   *_new -> *_FV2304
# Then wrap the data range in a native Excel Table ("Table1") and freeze
# the header row (matches the commit: "Use <formatversion> as suffix for
# table headers" / "Adjust xlsx export to new header formatting").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -----------------------------------
$suffixMap = @{
    "Segmentname_old"         = "Segmentname_FV2210"
    "Segmentgruppe_old"       = "Segmentgruppe_FV2210"
    "Segment_old"             = "Segment_FV2210"
    "Datenelement_old"        = "Datenelement_FV2210"
    "Segment ID_old"          = "Segment ID_FV2210"
    "Code_old"                = "Code_FV2210"
    "Qualifier_old"           = "Qualifier_FV2210"
    "Beschreibung_old"        = "Beschreibung_FV2210"
    "Bedingungsausdruck_old"  = "Bedingungsausdruck_FV2210"
    "Bedingung_old"           = "Bedingung_FV2210"
    "Segmentname_new"         = "Segmentname_FV2304"
    "Segmentgruppe_new"       = "Segmentgruppe_FV2304"
    "Segment_new"             = "Segment_FV2304"
    "Datenelement_new"        = "Datenelement_FV2304"
    "Segment ID_new"          = "Segment ID_FV2304"
    "Code_new"                = "Code_FV2304"
    "Qualifier_new"           = "Qualifier_FV2304"
    "Beschreibung_new"        = "Beschreibung_FV2304"
    "Bedingungsausdruck_new"  = "Bedingungsausdruck_FV2304"
    "Bedingung_new"           = "Bedingung_FV2304"
}

$lastRow = 82
$lastCol = 21

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Text
    if ($suffixMap.ContainsKey($old)) {
        $cell.Value = $suffixMap[$old]
    }
}

# --- 2. Turn the data range into a native Excel Table ----------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
